$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# About sheet: note the new state name ("Minnesota") and refresh the
# "last updated" date stamp in C1.
# ------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B1").Value = "Minnesota"
$wsAbout.Range("C1").Value = 44824

# ------------------------------------------------------------------
# Data sheet: the sheet only really uses columns A:D: columns E:J were
# leftover template formatting/placeholder cells with no real data.
# Clean them up - drop the stray formatting (and the handful of cells
# that held no content at all disappear entirely once cleared).
# ------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

# E3 holds a real value (20767) - only strip its formatting, keep the value.
$wsData.Range("E3").ClearFormats()

# The remaining E:J cells in these rows are empty placeholders - clear them
# completely so they drop out of the sheet.
$wsData.Range("F3:J3").Clear()
$wsData.Range("E4:J9").Clear()
$wsData.Range("E11:J11").Clear()
$wsData.Range("E13:J13").Clear()

# ------------------------------------------------------------------
# GDPGR-alternate: point the alternate growth-rate cell at the BAU
# sheet's computed rate instead of pulling straight from Data!B14.
# ------------------------------------------------------------------
$wsAlt = $wb.Worksheets.Item("GDPGR-alternate")
$wsAlt.Range("B2").Formula = "='GDPGR-bau'!B2"

# ------------------------------------------------------------------
# Tab selection moves from Data to GDPGR-alternate.
# ------------------------------------------------------------------
$wsAlt.Activate()
$wsAlt.Range("B3").Select()
